$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text could be misread as a number by Excel;
# force text format before assigning so they remain strings.
$numericLooking = @{
    "D4" = "1.008"
    "D5" = "334.15"
    "D7" = "0.4619"
    "D8" = "0.3859"
    "D9" = "45.98"
    "D10" = "0.07902"
    "D11" = "0.9965"
    "D12" = "21.48"
    "D13" = "5.961"
    "D15" = "7.126"
    "D18" = "0.06663"
    "D19" = "0.00001034"
    "D20" = "17.09"
    "D21" = "1.005"
    "D23" = "5.390"
    "D24" = "10.87"
    "D25" = "2.314"
    "D26" = "158.86"
    "D28" = "19.48"
    "D29" = "2.105"
    "D30" = "5.401"
    "D31" = "120.00"
    "D32" = "0.9731"
    "D33" = "0.09391"
    "D34" = "3.596"
    "D35" = "5.296"
    "D36" = "1.329"
    "D37" = "0.06014"
    "D39" = "8.257"
    "D40" = "1.180"
    "D41" = "0.5894"
    "D42" = "10.33"
    "D43" = "0.1857"
    "D45" = "0.5576"
    "D46" = "12.11"
    "D47" = "1.906"
    "D49" = "109.98"
}
foreach ($addr in $numericLooking.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLooking[$addr]
}

# Remaining text updates (coin names, links, prices, % change)
$textUpdates = @{
    "D2" = "27.500.11"
    "E2" = "  -1.50%  "
    "D3" = "1.842.05"
    "E3" = "  -2.05%  "
    "E4" = "  -1.19%  "
    "E5" = "  -0.20%  "
    "E6" = "  -1.22%  "
    "E8" = "  -0.93%  "
    "E9" = "  -1.88%  "
    "E10" = "  -0.13%  "
    "E11" = "  -0.72%  "
    "E12" = "  -0.14%  "
    "B13" = "Polkadot"
    "C13" = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
    "E13" = "  +0.96%  "
    "B14" = "WrappedEther"
    "C14" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "D14" = "1.833.38"
    "E14" = "  -3.37%  "
    "E15" = "  +1.00%  "
    "E16" = "  -1.37%  "
    "E17" = "  +1.89%  "
    "E18" = "  -1.18%  "
    "E19" = "  -0.43%  "
    "E20" = "  +0.79%  "
    "E21" = "  -1.21%  "
    "D22" = "27.495.12"
    "E22" = "  -1.60%  "
    "E23" = "  -0.98%  "
    "E24" = "  +0.13%  "
    "E25" = "  -1.54%  "
    "E26" = "  -0.38%  "
    "D27" = "2.063.84"
    "E27" = "  -2.67%  "
    "E28" = "  -1.71%  "
    "E29" = "  +2.69%  "
    "E31" = "  -0.59%  "
    "E32" = "  +2.13%  "
    "E33" = "  -0.64%  "
    "E34" = "  -1.80%  "
    "E35" = "  +0.18%  "
    "E36" = "  -1.30%  "
    "E37" = "  -1.16%  "
    "E38" = "  +0.01%  "
    "E39" = "  +2.52%  "
    "E40" = "  -1.87%  "
    "E41" = "  +0.68%  "
    "E42" = "  +2.44%  "
    "E43" = "  -0.88%  "
    "E44" = "  -2.34%  "
    "E46" = "  +0.74%  "
    "E47" = "  +0.26%  "
    "E48" = "  -2.93%  "
    "E49" = "  -2.78%  "
    "E50" = "  -1.06%  "
    "E51" = "  -1.33%  "
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}
